$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 72, shifting existing rows 72..194 down to 73..195
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new data record
$ws.Range("A72").Value = 4
$ws.Range("B72").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C72").Value = "Los Lagos"
$ws.Range("D72").Value = 44477
$ws.Range("E72").Value = 10
$ws.Range("F72").Value = 100112008
$ws.Range("G72").Value = "Coliflor"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 1200
$ws.Range("K72").Value = 1000
$ws.Range("L72").Value = 1100
$ws.Range("M72").Value = 1050
$ws.Range("N72").Value = '$/unidad'
$ws.Range("O72").Value = "Región Metropolitana"
$ws.Range("P72").Value = 1050
$ws.Range("Q72").Value = 1
$ws.Range("R72").Value = "Hortaliza"

Write-Host "Done"
